$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61 - this shifts existing rows 61:85 down to 62:86
$ws.Rows.Item(61).Insert()

# Populate the newly inserted row 61 with the new weekly record
$ws.Cells.Item(61, 1).Value = 10
$ws.Cells.Item(61, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(61, 3).Value = "La Araucanía"
$ws.Cells.Item(61, 4).NumberFormat = $ws.Cells.Item(62, 4).NumberFormat
$ws.Cells.Item(61, 4).Value = 44726
$ws.Cells.Item(61, 5).Value = 9
$ws.Cells.Item(61, 6).Value = 100112035
$ws.Cells.Item(61, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 35
$ws.Cells.Item(61, 11).Value = 30000
$ws.Cells.Item(61, 12).Value = 30000
$ws.Cells.Item(61, 13).Value = 30000
$ws.Cells.Item(61, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(61, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(61, 16).Value = 3000
$ws.Cells.Item(61, 17).Value = 10
$ws.Cells.Item(61, 18).Value = "Hortaliza"
